$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 3901.49
$ws.Range("E2").Value = -3901.49

$ws.Range("D4").Value = 4449.95
$ws.Range("E4").Value = 13050.05
$ws.Range("F4").Value = 0.2542828571428571
